$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column updates)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2609
$ws.Range("F3").Value = 565
$ws.Range("F4").Value = 465
$ws.Range("F6").Value = 190
$ws.Range("F7").Value = 470
$ws.Range("F8").Value = 1210
$ws.Range("F9").Value = 552
$ws.Range("F10").Value = 301
$ws.Range("F11").Value = 121
$ws.Range("F12").Value = 355
$ws.Range("F13").Value = 5643
$ws.Range("F14").Value = 72
$ws.Range("F15").Value = 1735
$ws.Range("F16").Value = 4085
$ws.Range("F17").Value = 421
$ws.Range("F18").Value = 238
$ws.Range("F19").Value = 303
$ws.Range("F20").Value = 4717
$ws.Range("F21").Value = 6138
$ws.Range("F23").Value = 1043
$ws.Range("F24").Value = 675
$ws.Range("F25").Value = 3729
$ws.Range("F26").Value = 490
$ws.Range("F27").Value = 67
$ws.Range("F28").Value = 187
$ws.Range("F29").Value = 123
$ws.Range("F30").Value = 978
$ws.Range("F31").Value = 1387
$ws.Range("F32").Value = 459
$ws.Range("F33").Value = 537
$ws.Range("F34").Value = 1579
$ws.Range("F36").Value = 1687
$ws.Range("F37").Value = 179
$ws.Range("F39").Value = 1110
$ws.Range("F40").Value = 32
$ws.Range("F41").Value = 1340
$ws.Range("F42").Value = 621
$ws.Range("F44").Value = 3352
$ws.Range("F46").Value = 274
$ws.Range("F48").Value = 6
$ws.Range("F49").Value = 3874

# Sheet "演出" (F column updates)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1190
$ws.Range("F7").Value = 7
$ws.Range("F10").Value = 14

# Sheet "本地生活" (F column updates)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3792

# Sheet "全部类型" (F column updates)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3792
$ws.Range("F4").Value = 2609
$ws.Range("F5").Value = 565
$ws.Range("F6").Value = 465
$ws.Range("F8").Value = 1190
$ws.Range("F9").Value = 7
$ws.Range("F10").Value = 190
$ws.Range("F11").Value = 470
$ws.Range("F12").Value = 1210
$ws.Range("F13").Value = 552
$ws.Range("F14").Value = 301
$ws.Range("F15").Value = 121
$ws.Range("F16").Value = 355
$ws.Range("F18").Value = 1735
$ws.Range("F19").Value = 4717
$ws.Range("F21").Value = 1043
$ws.Range("F22").Value = 675
$ws.Range("F23").Value = 3729
$ws.Range("F24").Value = 490
$ws.Range("F25").Value = 67
$ws.Range("F26").Value = 187
$ws.Range("F27").Value = 123
$ws.Range("F28").Value = 978
$ws.Range("F29").Value = 1387
$ws.Range("F30").Value = 459
$ws.Range("F31").Value = 537
$ws.Range("F33").Value = 1579
$ws.Range("F35").Value = 1687
$ws.Range("F37").Value = 1110
$ws.Range("F39").Value = 621
$ws.Range("F43").Value = 3352
$ws.Range("F46").Value = 274
$ws.Range("F49").Value = 3874
